$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.691.21"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "1.685.38"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'220.12"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'30.39"
$ws.Range("D9").Value = "'0.265"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "1.928.19"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("D13").Value = "'10.45"
$ws.Range("E13").Value = "  +11.39%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.619"
$ws.Range("E14").Value = "  +8.15%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.686.47"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "30.698.76"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "'66.26"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "'246.06"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'157.66"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'15.83"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "'3.47"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "1.509.62"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "'3.29"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "'84.50"
$ws.Range("E36").Value = "  +9.20%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "'0.584"
$ws.Range("E40").Value = "  +4.25%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'0.836"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'51.86"
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("D48").Value = "1.820.06"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'94.74"
$ws.Range("E50").Value = "  +5.05%  "
$ws.Range("D51").Value = "0.0₆0115"
$ws.Range("E51").Value = "  +0.70%  "
